# Calibration of energy use modeling by renovation level
# Apply a 1.12x calibration factor to the "AB" UFA values for rows 403-452
# (years 2001-2050), matching the committed change.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$factor = 1.12

for ($row = 403; $row -le 452; $row++) {
    $cell = $ws.Cells.Item($row, 2)
    $current = $cell.Value2
    $cell.Value = $current * $factor
}
